$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.437724
$ws.Range("H2").Value = 1.313172
$ws.Range("I2").Value = 0.02046276855287852
$ws.Range("J2").Value = 0.02204588088728605
$ws.Range("Q2").Value = 0.01652991732
$ws.Range("R2").Value = 0.14876925588
$ws.Range("S2").Value = 0.02046276855287852
$ws.Range("T2").Value = 0.02204588088728605

# Row 3
$ws.Range("G3").Value = 12.48419333333333
$ws.Range("I3").Value = 0.5836124104444559
$ws.Range("J3").Value = 0.6287638767819841
$ws.Range("Q3").Value = 0.4714447542444444
$ws.Range("R3").Value = 4.2430027882
$ws.Range("S3").Value = 0.5836124104444559
$ws.Range("T3").Value = 0.6287638767819841

# Row 4
$ws.Range("G4").Value = 1.796802333333333
$ws.Range("H4").Value = 5.390407
$ws.Range("I4").Value = 0.08399710841140098
$ws.Range("J4").Value = 0.09049558675938332
$ws.Range("Q4").Value = 0.06785324544777778
$ws.Range("R4").Value = 0.6106792090299999
$ws.Range("S4").Value = 0.08399710841140098
$ws.Range("T4").Value = 0.09049558675938332

# Row 5
$ws.Range("G5").Value = 4.608308
$ws.Range("H5").Value = 9.216616
$ws.Range("I5").Value = 0.2154296772038511
$ws.Range("J5").Value = 0.154731001361478
$ws.Range("Q5").Value = 0.1740250711066667
$ws.Range("R5").Value = 1.04415042664
$ws.Range("S5").Value = 0.2154296772038511
$ws.Range("T5").Value = 0.154731001361478

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 2.064212666666667
$ws.Range("H6").Value = 6.192638000000001
$ws.Range("I6").Value = 0.09649803538741349
$ws.Range("J6").Value = 0.1039636542098684
$ws.Range("Q6").Value = 0.07795155100222223
$ws.Range("R6").Value = 0.70156395902
$ws.Range("S6").Value = 0.09649803538741349
$ws.Range("T6").Value = 0.1039636542098684

$wb.Save()
